# Apply GitHub-Actions-style crypto price/volume refresh to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) cells hold free-form text (e.g. "29.818.90", "0.990") rather than
# numeric values, so force Text format before writing to avoid Excel re-interpreting
# the strings as numbers (which would drop significant trailing zeros).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = '29.818.90'
$ws.Range("E2").Value = '  +1.13%  '
$ws.Range("D3").Value = '1.619.22'
$ws.Range("E3").Value = '  +0.77%  '
$ws.Range("E4").Value = '  -0.85%  '
$ws.Range("D5").Value = '213.06'
$ws.Range("E6").Value = '  -0.24%  '
$ws.Range("D7").Value = '0.990'
$ws.Range("D8").Value = '29.14'
$ws.Range("E8").Value = '  +8.44%  '
$ws.Range("D9").Value = '0.258'
$ws.Range("E9").Value = '  +2.91%  '
$ws.Range("E10").Value = '  +0.89%  '
$ws.Range("E11").Value = '  -0.05%  '
$ws.Range("D12").Value = '1.852.06'
$ws.Range("D13").Value = '1.608.53'
$ws.Range("E13").Value = '  -0.42%  '
$ws.Range("E14").Value = '  +5.61%  '
$ws.Range("E15").Value = '  +4.87%  '
$ws.Range("D16").Value = '29.827.89'
$ws.Range("E16").Value = '  +1.07%  '
$ws.Range("D17").Value = '8.89'
$ws.Range("E17").Value = '  +16.67%  '
$ws.Range("D18").Value = '64.41'
$ws.Range("E18").Value = '  +1.81%  '
$ws.Range("D19").Value = '241.36'
$ws.Range("E19").Value = '  +0.07%  '
$ws.Range("E20").Value = '  +2.40%  '
$ws.Range("D21").Value = '0.993'
$ws.Range("E21").Value = '  -0.59%  '
$ws.Range("D22").Value = '4.10'
$ws.Range("E22").Value = '  +2.47%  '
$ws.Range("D23").Value = '9.60'
$ws.Range("E23").Value = '  +4.39%  '
$ws.Range("D24").Value = '2.11'
$ws.Range("E24").Value = '  +0.67%  '
$ws.Range("D25").Value = '154.91'
$ws.Range("E25").Value = '  +0.31%  '
$ws.Range("D26").Value = '15.61'
$ws.Range("E26").Value = '  +2.26%  '
$ws.Range("D27").Value = '0.110'
$ws.Range("E27").Value = '  +1.27%  '
$ws.Range("E28").Value = '  +3.30%  '
$ws.Range("E29").Value = '  -0.74%  '
$ws.Range("E30").Value = '  +2.99%  '
$ws.Range("D31").Value = '1.10'
$ws.Range("E31").Value = '  +3.59%  '
$ws.Range("D32").Value = '3.33'
$ws.Range("E32").Value = '  +3.39%  '
$ws.Range("D33").Value = '3.21'
$ws.Range("E33").Value = '  +3.57%  '
$ws.Range("D34").Value = '1.416.15'
$ws.Range("E34").Value = '  +0.27%  '
$ws.Range("D35").Value = '1.64'
$ws.Range("E35").Value = '  +6.75%  '
$ws.Range("E36").Value = '  +0.05%  '
$ws.Range("D37").Value = '2.87'
$ws.Range("E37").Value = '  +1.19%  '
$ws.Range("D38").Value = '2.29'
$ws.Range("E38").Value = '  -0.65%  '
$ws.Range("E39").Value = '  +2.54%  '
$ws.Range("D40").Value = '0.556'
$ws.Range("E40").Value = '  +3.39%  '
$ws.Range("D41").Value = '0.0503'
$ws.Range("E41").Value = '  +3.59%  '
$ws.Range("D42").Value = '1.98'
$ws.Range("E42").Value = '  +0.16%  '
$ws.Range("D44").Value = '53.89'
$ws.Range("E44").Value = '  +2.24%  '
$ws.Range("D45").Value = '69.35'
$ws.Range("E45").Value = '  +5.67%  '
$ws.Range("E46").Value = '  +18.48%  '
$ws.Range("D47").Value = '0.991'
$ws.Range("E47").Value = '  -0.73%  '
$ws.Range("D48").Value = '5.43'
$ws.Range("E48").Value = '  +2.90%  '
$ws.Range("D49").Value = '1.759.81'
$ws.Range("E49").Value = '  +0.60%  '
$ws.Range("D50").Value = '88.12'
$ws.Range("E50").Value = '  +1.65%  '
$ws.Range("E51").Value = '  +2.04%  '
